$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.88"
$ws.Range("E2").Value = "'5.12%"
$ws.Range("D3").Value = "'31.92"
$ws.Range("E3").Value = "'9.38%"
$ws.Range("D4").Value = "'5.271"
$ws.Range("E4").Value = "'0.06%"
$ws.Range("D5").Value = "'0.07488"
$ws.Range("E5").Value = "'7.06%"
$ws.Range("E6").Value = "'5.35%"
$ws.Range("D7").Value = "'3.816"
$ws.Range("E7").Value = "'7.41%"
$ws.Range("D8").Value = "'1.488"
$ws.Range("E8").Value = "'7.07%"
$ws.Range("D9").Value = "'0.9201"
$ws.Range("E9").Value = "'1.78%"
$ws.Range("D10").Value = "'0.1684"
$ws.Range("E10").Value = "'5.07%"
$ws.Range("D11").Value = "'0.07870"
$ws.Range("E11").Value = "'3.78%"
$ws.Range("D12").Value = "'0.08035"
$ws.Range("E12").Value = "'4.02%"
$ws.Range("E13").Value = "'3.65%"
$ws.Range("D14").Value = "'0.09882"
$ws.Range("E14").Value = "'9.42%"
$ws.Range("D15").Value = "'0.001504"
$ws.Range("E15").Value = "'-5.43%"
$ws.Range("D16").Value = "'0.04604"
$ws.Range("E16").Value = "'1.74%"
$ws.Range("D17").Value = "'0.006367"
$ws.Range("E17").Value = "'-0.86%"
$ws.Range("D18").Value = "'3.460"
$ws.Range("E18").Value = "'-0.87%"
$ws.Range("E19").Value = "'-0.12%"
$ws.Range("D20").Value = "'0.3301"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("E21").Value = "'-0.01%"
$ws.Range("D22").Value = "'4.499"
$ws.Range("E22").Value = "'12.04%"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("D25").Value = "'0.004442"
$ws.Range("E25").Value = "'7.11%"
$ws.Range("D26").Value = "'0.0001399"
$ws.Range("E26").Value = "'19.66%"
$ws.Range("D27").Value = "'0.0001774"
$ws.Range("E27").Value = "'6.38%"
$ws.Range("D39").Value = "'0.01714"
$ws.Range("E39").Value = "'2,528.65%"
$ws.Range("D40").Value = "'0.04480"
$ws.Range("E40").Value = "'2.69%"
$ws.Range("D41").Value = "'0.006982"
$ws.Range("E41").Value = "'0.74%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'7.78%"
$ws.Range("D43").Value = "'0.002249"
$ws.Range("E43").Value = "'8.71%"
$ws.Range("D44").Value = "'0.01282"
$ws.Range("E44").Value = "'10.30%"
$ws.Range("D45").Value = "'0.00006154"
$ws.Range("E45").Value = "'5.66%"
$ws.Range("D46").Value = "'0.7111"
$ws.Range("E46").Value = "'-63.14%"
$ws.Range("D47").Value = "'0.01498"
$ws.Range("E47").Value = "'15.21%"
